$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.821.30"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "1.869.90"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'301.02"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.5341"
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("D8").Value = "'0.3749"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "'0.07172"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "'21.53"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "'0.8861"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "'0.08145"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "1.877.92"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "'92.97"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "'5.250"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "'14.69"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'0.000008534"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "26.889.94"
$ws.Range("D21").Value = "'4.967"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").Value = "'10.69"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'6.389"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "'146.96"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "'2.254"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.00"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.725"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "'114.24"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "'4.740"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").Value = "'4.571"
$ws.Range("E30").Value = "  -6.66%  "
$ws.Range("D31").Value = "'0.09127"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "'0.7987"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").Value = "'2.993"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "'1.169"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("D36").Value = "'0.5986"
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("D37").Value = "'2.604"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'3.146"
$ws.Range("E38").Value = "  -6.43%  "
$ws.Range("D39").Value = "'0.01946"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "'1.070"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "'6.638"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").Value = "'8.880"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "'0.5101"
$ws.Range("E44").Value = "  +4.42%  "
$ws.Range("D45").Value = "'0.1495"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'9.941"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "'1.625"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "'37.67"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").Value = "'0.06016"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'62.14"
